{"js": "// Wrap the \"Choix d'un fichier de param\u00e8tre de connexion \u00e0 Discord\" bullet\n// text in parentheses and move the \"_GoBack\" bookmark from the end of the\n// \"Serveur (ChatBot) :\" paragraph to just after the new opening \"(\" run.\n\nconst TARGET_TEXT = \"Choix d\\u2019un fichier de param\\u00e8tre de connexion \\u00e0 Discord\";\nconst BOOKMARK_NAME = \"_GoBack\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(TARGET_TEXT) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the target paragraph: \" + TARGET_TEXT);\n}\n\n// Remove the existing \"_GoBack\" bookmark first so re-inserting it below\n// leaves exactly one bookmark of that name (moved, not duplicated).\ncontext.document.deleteBookmark(BOOKMARK_NAME);\nawait context.sync();\n\n// Insert the opening \"(\" as its own run right before the paragraph's text.\nconst openParenRange = target.getRange(\"Start\").insertText(\"(\", \"Before\");\nawait context.sync();\n\n// Re-create the \"_GoBack\" bookmark immediately after the \"(\" run (collapsed,\n// i.e. bookmarkStart immediately followed by bookmarkEnd).\nconst afterOpenParen = openParenRange.getRange(\"After\");\nafterOpenParen.insertBookmark(BOOKMARK_NAME);\nawait context.sync();\n\n// Append the closing \")\" as its own run at the end of the paragraph.\ntarget.getRange(\"End\").insertText(\")\", \"After\");\nawait context.sync();\n", "ps1": "# Wrap the \"Choix d'un fichier de param\u00e8tre de connexion \u00e0 Discord\" bullet\n# text in parentheses and move the \"_GoBack\" bookmark from the end of the\n# \"Serveur (ChatBot) :\" paragraph to just after the new opening \"(\" run.\n\n$d = $word.ActiveDocument\n$bookmarkName = \"_GoBack\"\n\n# 1. Find the target paragraph containing \"...connexion \u00e0 Discord\".\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Discord*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the target paragraph containing 'Discord'\"\n}\n\n# 2. Remove the existing \"_GoBack\" bookmark first so re-adding it below\n#    leaves exactly one bookmark of that name (moved, not duplicated).\nif ($d.Bookmarks.Exists($bookmarkName)) {\n    $d.Bookmarks.Item($bookmarkName).Delete()\n}\n\n# 3. Insert the opening \"(\" as its own run right before the paragraph's text.\n$target.Range.InsertBefore(\"(\")\n\n# 4. Re-create the \"_GoBack\" bookmark immediately after the \"(\" run\n#    (collapsed range, i.e. bookmarkStart immediately followed by bookmarkEnd).\n$pStart = $target.Range.Start\n$bmRange = $d.Range($pStart + 1, $pStart + 1)\n$d.Bookmarks.Add($bookmarkName, $bmRange)\n\n# 5. Append the closing \")\" as its own run at the end of the paragraph.\n#    (End - 1 lands right before the paragraph mark, i.e. still inside\n#    this paragraph rather than spilling into the next one.)\n$pEnd = $target.Range.End - 1\n$endRange = $d.Range($pEnd, $pEnd)\n$endRange.InsertAfter(\")\")\n"}
